# Add data for 2022-10-20
# - Rename the "through" date from October 11 to October 12 (sheet title + label)
# - Bump several existing neighborhood/month counts
# - Add a handful of brand new neighborhood/month counts (previously empty cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Sheet name and the matching column header / shared string label
$ws.Name = "Through 2022-10-12"
$ws.Range("B1").Value = "October 2022 (through October 12)"

# 2. Existing cell value bumps
$ws.Range("L2").Value = 9    # Garfield Park       - October 2021
$ws.Range("V2").Value = 7    # Garfield Park       - October 2020
$ws.Range("B3").Value = 3    # Englewood           - October 2022 (through October 12)
$ws.Range("V3").Value = 4    # Englewood           - October 2020
$ws.Range("L4").Value = 5    # Humboldt Park       - October 2021
$ws.Range("L5").Value = 6    # Austin              - October 2021
$ws.Range("V6").Value = 6    # North Lawndale      - October 2020
$ws.Range("V7").Value = 2    # South Shore         - October 2020
$ws.Range("AP7").Value = 2   # South Shore         - October 2018
$ws.Range("L29").Value = 2   # Lake View           - October 2021
$ws.Range("L52").Value = 2   # Irving Park         - October 2021

# 3. Brand-new counts (cells that were previously blank)
$ws.Range("BT10").Value = 1  # Douglas             - October 2015
$ws.Range("AF18").Value = 1  # Washington Heights  - October 2019
$ws.Range("AZ23").Value = 1  # Auburn Gresham      - October 2017
$ws.Range("BT23").Value = 1  # Auburn Gresham      - October 2015
$ws.Range("AF29").Value = 1  # Lake View           - October 2019
$ws.Range("L39").Value = 1   # Albany Park         - October 2021
$ws.Range("B52").Value = 1   # Irving Park         - October 2022 (through October 12)
$ws.Range("AF96").Value = 1  # West Elsdon         - October 2019
